$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add three new header cells, matching the existing header style ---
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"

$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: new data row ---
$ws.Range("A2").Value = "parisk"
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = "APC"
$ws.Range("E2").Value = "WRI"
$ws.Range("F2").Value = "1269f1fb-9c21-42a9-ae5e-c80f92622adc"
$ws.Range("G2").Value = "Bk6qQGWRb_annotated.xlsx"
$ws.Range("H2").Value = "Then how bootstrap dqn extend the idea to deep learning, followed by the noisy net, bbq, shallow UBE and LS-DQN."

# C2 is an empty (but present) text cell in the source data. A bare leading
# apostrophe creates an empty text-typed cell; copy the plain (unstyled)
# format from a sibling cell so C2 doesn't pick up a quote-prefix style.
$ws.Range("C2").Value = "'"
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
